$d = $word.ActiveDocument

# The document currently has a single paragraph:
#   "Test" + bookmarkStart/_GoBack + bookmarkEnd
#
# We need to split it into three paragraphs:
#   1) "Test"
#   2) "Adding Test Security Data"
#   3) (empty, holding the _GoBack bookmark)
#
# Using Find & Replace with "^p" paragraph-mark codes lets the new
# paragraph marks be inserted right after "Test", pushing the existing
# bookmark down into its own trailing paragraph, and creating a fresh
# paragraph in between that carries the new sentence.
$d.Content.Find.Execute(
    "Test",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Test^pAdding Test Security Data^p",
    2
)
